$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure the Price column (D) is treated as text so values like "1.00" or
# "0.0000127" keep their exact formatting instead of being parsed as numbers.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "66.677.77"
$ws.Range("E2").Value = "  +0.75%  "
$ws.Range("D3").Value = "3.454.13"
$ws.Range("E3").Value = "  +0.48%  "
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("D5").Value = "584.99"
$ws.Range("E5").Value = "  +0.93%  "
$ws.Range("D6").Value = "178.09"
$ws.Range("E6").Value = "  +2.92%  "
$ws.Range("D7").Value = "0.630"
$ws.Range("E7").Value = "  +6.92%  "
$ws.Range("E8").Value = "  +0.00%  "
$ws.Range("D9").Value = "3.455.00"
$ws.Range("E9").Value = "  +0.63%  "
$ws.Range("E10").Value = "  +1.59%  "
$ws.Range("E11").Value = "  +1.47%  "
$ws.Range("D12").Value = "0.418"
$ws.Range("E12").Value = "  -0.06%  "
$ws.Range("D13").Value = "4.062.87"
$ws.Range("E13").Value = "  +0.54%  "
$ws.Range("E14").Value = "  +1.77%  "
$ws.Range("D15").Value = "30.18"
$ws.Range("E15").Value = "  +0.60%  "
$ws.Range("D16").Value = "66.563.45"
$ws.Range("E16").Value = "  +0.53%  "
$ws.Range("E17").Value = "  +1.76%  "
$ws.Range("D18").Value = "3.455.11"
$ws.Range("E18").Value = "  -0.05%  "
$ws.Range("D19").Value = "5.97"
$ws.Range("E19").Value = "  +0.02%  "
$ws.Range("D20").Value = "13.89"
$ws.Range("E20").Value = "  +1.46%  "
$ws.Range("D21").Value = "372.02"
$ws.Range("E21").Value = "  -1.43%  "
$ws.Range("D22").Value = "7.66"
$ws.Range("E22").Value = "  -0.99%  "
$ws.Range("D23").Value = "73.43"
$ws.Range("E23").Value = "  +2.09%  "
$ws.Range("D24").Value = "0.0000127"
$ws.Range("E24").Value = "  +8.21%  "
$ws.Range("D25").Value = "0.998"
$ws.Range("E25").Value = "  -0.02%  "
$ws.Range("D26").Value = "0.535"
$ws.Range("E26").Value = "  -1.41%  "
$ws.Range("D27").Value = "9.94"
$ws.Range("E27").Value = "  +1.82%  "
$ws.Range("E28").Value = "  +2.86%  "
$ws.Range("D29").Value = "1.00"
$ws.Range("E29").Value = "  +0.20%  "
$ws.Range("D30").Value = "5.90"
$ws.Range("E30").Value = "  +1.46%  "
$ws.Range("E31").Value = "  +0.84%  "
$ws.Range("D32").Value = "23.64"
$ws.Range("E32").Value = "  -1.85%  "
$ws.Range("D33").Value = "0.999"
$ws.Range("E33").Value = "  -0.05%  "
$ws.Range("E34").Value = "  -0.15%  "
$ws.Range("D35").Value = "1.27"
$ws.Range("E35").Value = "  -2.42%  "
$ws.Range("E36").Value = "  +0.74%  "
$ws.Range("D37").Value = "162.94"
$ws.Range("E37").Value = "  +2.34%  "
$ws.Range("D38").Value = "0.883"
$ws.Range("E38").Value = "  +0.22%  "
$ws.Range("D39").Value = "27.85"
$ws.Range("E39").Value = "  -3.90%  "
$ws.Range("E40").Value = "  +2.77%  "
$ws.Range("E41").Value = "  +3.38%  "
$ws.Range("D42").Value = "4.50"
$ws.Range("E42").Value = "  +0.53%  "
$ws.Range("D43").Value = "2.759.77"
$ws.Range("E43").Value = "  +3.84%  "
$ws.Range("D44").Value = "6.48"
$ws.Range("E44").Value = "  +1.98%  "
$ws.Range("D45").Value = "0.0696"
$ws.Range("E45").Value = "  +0.92%  "
$ws.Range("D46").Value = "25.40"
$ws.Range("E46").Value = "  +4.91%  "
$ws.Range("D47").Value = "340.97"
$ws.Range("E47").Value = "  +9.65%  "
$ws.Range("D48").Value = "40.06"
$ws.Range("E48").Value = "  -0.52%  "
$ws.Range("D49").Value = "0.0288"
$ws.Range("E49").Value = "  -0.12%  "
$ws.Range("D50").Value = "0.105"
$ws.Range("D51").Value = "31.77"
$ws.Range("E51").Value = "  +4.41%  "
